$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 75
$ws.Cells.Item($row, 1).Value = "2024-10-16 00:00:00"
$ws.Cells.Item($row, 2).Value = 75400
$ws.Cells.Item($row, 3).Value = 10561.7
$ws.Cells.Item($row, 4).Value = 9346.639999999999
$ws.Cells.Item($row, 5).Value = 7.1125
